$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.538652062416077
$ws.Range("B1").Value = 3.596222162246704
$ws.Range("C1").Value = 4.806467056274414
$ws.Range("D1").Value = 1.986483931541443
$ws.Range("E1").Value = 1.077107310295105
